$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Add a new user row (row 52) with a new "Linking_AutoUser" test account,
# mirroring the existing data rows (e.g. row 51).
$row = 52
$ws.Cells.Item($row, 1).Value = "Linking_AutoUser"
$ws.Cells.Item($row, 2).Value = "Password1"
$ws.Cells.Item($row, 5).Value = "Default user for Linking tests"
$ws.Cells.Item($row, 6).Value = "N"
$ws.Cells.Item($row, 7).Value = "linking.autouser@mailinator.com"

# Match the bordered look of the rest of the table.
$rng = $ws.Range("A52:G52")
$rng.Borders.Color = 0
$rng.Borders.LineStyle = 1

# Reflect the new selection/scroll position recorded for the sheet.
$ws.Range("C26").Select() | Out-Null
